$wb = $excel.ActiveWorkbook

# ---- Sheet1: summary table ----
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Item(2, 2).Value = -411.5189726818494
$ws1.Cells.Item(2, 3).Value = 15.515153592
$ws1.Cells.Item(2, 4).Value = 3
$ws1.Cells.Item(2, 6).Value = 2
$ws1.Cells.Item(2, 7).Value = 1204
$ws1.Cells.Item(2, 8).Value = 1210
$ws1.Cells.Item(2, 9).Value = 100

$ws1.Cells.Item(3, 2).Value = -411.39604929747
$ws1.Cells.Item(3, 3).Value = 2.087831791
$ws1.Cells.Item(3, 4).Value = 4
$ws1.Cells.Item(3, 6).Value = 3
$ws1.Cells.Item(3, 7).Value = 1756
$ws1.Cells.Item(3, 8).Value = 1815
$ws1.Cells.Item(3, 9).Value = 150

$ws1.Cells.Item(4, 2).Value = -412.34762188775494
$ws1.Cells.Item(4, 3).Value = 3.424428263
$ws1.Cells.Item(4, 4).Value = 5
$ws1.Cells.Item(4, 6).Value = 4
$ws1.Cells.Item(4, 7).Value = 2308
$ws1.Cells.Item(4, 8).Value = 2420
$ws1.Cells.Item(4, 9).Value = 200

$ws1.Cells.Item(5, 2).Value = -416.71129461891667
$ws1.Cells.Item(5, 3).Value = 12.33823612
$ws1.Cells.Item(5, 4).Value = 4
$ws1.Cells.Item(5, 6).Value = 3
$ws1.Cells.Item(5, 7).Value = 1756
$ws1.Cells.Item(5, 8).Value = 1815
$ws1.Cells.Item(5, 9).Value = 150

$ws1.Cells.Item(6, 2).Value = -408.2905898186462
$ws1.Cells.Item(6, 3).Value = 2.134491656
$ws1.Cells.Item(6, 4).Value = 4
$ws1.Cells.Item(6, 6).Value = 3
$ws1.Cells.Item(6, 7).Value = 1756
$ws1.Cells.Item(6, 8).Value = 1815
$ws1.Cells.Item(6, 9).Value = 150

$ws1.Cells.Item(7, 2).Value = -403.9902935908759
$ws1.Cells.Item(7, 3).Value = 3.447364775
$ws1.Cells.Item(7, 4).Value = 4
$ws1.Cells.Item(7, 6).Value = 3
$ws1.Cells.Item(7, 7).Value = 1756
$ws1.Cells.Item(7, 8).Value = 1815
$ws1.Cells.Item(7, 9).Value = 150

$ws1.Cells.Item(8, 2).Value = -400.1919964156508
$ws1.Cells.Item(8, 3).Value = 3.729217923
$ws1.Cells.Item(8, 4).Value = 2
$ws1.Cells.Item(8, 6).Value = 1
$ws1.Cells.Item(8, 7).Value = 652
$ws1.Cells.Item(8, 8).Value = 605
$ws1.Cells.Item(8, 9).Value = 50

$ws1.Cells.Item(9, 2).Value = -412.09265747622965
$ws1.Cells.Item(9, 3).Value = 9.611153647
$ws1.Cells.Item(9, 4).Value = 4
$ws1.Cells.Item(9, 6).Value = 3
$ws1.Cells.Item(9, 7).Value = 1756
$ws1.Cells.Item(9, 8).Value = 1815
$ws1.Cells.Item(9, 9).Value = 150

$ws1.Cells.Item(10, 2).Value = -408.4459826348807
$ws1.Cells.Item(10, 3).Value = 7.337151256
$ws1.Cells.Item(10, 4).Value = 7
$ws1.Cells.Item(10, 6).Value = 6
$ws1.Cells.Item(10, 7).Value = 3412
$ws1.Cells.Item(10, 8).Value = 3630
$ws1.Cells.Item(10, 9).Value = 300

$ws1.Cells.Item(11, 2).Value = -403.1039288610207
$ws1.Cells.Item(11, 3).Value = 16.283771156
$ws1.Cells.Item(11, 4).Value = 3
$ws1.Cells.Item(11, 6).Value = 2
$ws1.Cells.Item(11, 7).Value = 1204
$ws1.Cells.Item(11, 8).Value = 1210
$ws1.Cells.Item(11, 9).Value = 100

# ---- Sheet "1" detail ----
$wsX = $wb.Worksheets.Item("1")
$wsX.Cells.Item(2, 4).Value = 0.8748163855462646
$wsX.Cells.Item(2, 5).Value = 84.92934
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -411.8636030575577
$wsX.Cells.Item(3, 3).Value = 0.0
$wsX.Cells.Item(3, 4).Value = 1.1653051579102782
$wsX.Cells.Item(3, 5).Value = 1.18533
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -411.5189726818494
$wsX.Cells.Item(4, 3).Value = 0.0
$wsX.Cells.Item(4, 4).Value = 6.096441480468628
$wsX.Cells.Item(4, 5).Value = 0.0

# ---- Sheet "2" detail ----
$wsX = $wb.Worksheets.Item("2")
$wsX.Cells.Item(2, 4).Value = 0.03201625409387207
$wsX.Cells.Item(2, 5).Value = 81.97293
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -417.900310782136
$wsX.Cells.Item(3, 3).Value = 0.07686770769276234
$wsX.Cells.Item(3, 4).Value = 0.30459583754467773
$wsX.Cells.Item(3, 5).Value = 1.93166
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -415.4220225388203
$wsX.Cells.Item(4, 3).Value = 0.09923403939962108
$wsX.Cells.Item(4, 4).Value = 0.44836400724890135
$wsX.Cells.Item(4, 5).Value = 2.62995
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -411.39604929747
$wsX.Cells.Item(5, 3).Value = 0.08528937732092581
$wsX.Cells.Item(5, 4).Value = 1.1096881597039796
$wsX.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "3" detail ----
$wsX = $wb.Worksheets.Item("3")
$wsX.Cells.Item(2, 4).Value = 0.02876627896044922
$wsX.Cells.Item(2, 5).Value = 86.44331
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -420.50506127605115
$wsX.Cells.Item(3, 3).Value = 0.07256599914889564
$wsX.Cells.Item(3, 4).Value = 0.1901683605397949
$wsX.Cells.Item(3, 5).Value = 2.56965
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -416.31345173383806
$wsX.Cells.Item(4, 3).Value = 0.02562620199290218
$wsX.Cells.Item(4, 4).Value = 0.27396762162390137
$wsX.Cells.Item(4, 5).Value = 1.41192
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -412.41404327794714
$wsX.Cells.Item(5, 3).Value = 0.07066505111364295
$wsX.Cells.Item(5, 4).Value = 0.845777416161621
$wsX.Cells.Item(5, 5).Value = 0.19822
$wsX.Cells.Item(6, 1).Value = 5
$wsX.Cells.Item(6, 2).Value = -412.34762188775494
$wsX.Cells.Item(6, 3).Value = 0.01610810555678573
$wsX.Cells.Item(6, 4).Value = 1.7630294423380126
$wsX.Cells.Item(6, 5).Value = 0.0

# ---- Sheet "4" detail ----
$wsX = $wb.Worksheets.Item("4")
$wsX.Cells.Item(2, 4).Value = 0.01973398877355957
$wsX.Cells.Item(2, 5).Value = 85.55874
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -417.96361691595314
$wsX.Cells.Item(3, 3).Value = 0.08153552595003252
$wsX.Cells.Item(3, 4).Value = 2.2727969560262453
$wsX.Cells.Item(3, 5).Value = 2.36443
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -416.89458935220637
$wsX.Cells.Item(4, 3).Value = 0.002793873737720017
$wsX.Cells.Item(4, 4).Value = 3.208835709969849
$wsX.Cells.Item(4, 5).Value = 1.00555
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -416.71129461891667
$wsX.Cells.Item(5, 3).Value = 0.04678112742630615
$wsX.Cells.Item(5, 4).Value = 6.601189085973511
$wsX.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "5" detail ----
$wsX = $wb.Worksheets.Item("5")
$wsX.Cells.Item(2, 4).Value = 0.04153228308178711
$wsX.Cells.Item(2, 5).Value = 84.60907
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -448.74515541187657
$wsX.Cells.Item(3, 3).Value = 0.03707763961206097
$wsX.Cells.Item(3, 4).Value = 0.1310620778588867
$wsX.Cells.Item(3, 5).Value = 5.70689
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -409.5100971023726
$wsX.Cells.Item(4, 3).Value = 0.09813048616550464
$wsX.Cells.Item(4, 4).Value = 0.5924801040747071
$wsX.Cells.Item(4, 5).Value = 1.80117
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -408.2905898186462
$wsX.Cells.Item(5, 3).Value = 0.09134710180965602
$wsX.Cells.Item(5, 4).Value = 1.13146167544458
$wsX.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "6" detail ----
$wsX = $wb.Worksheets.Item("6")
$wsX.Cells.Item(2, 4).Value = 0.03483292279101562
$wsX.Cells.Item(2, 5).Value = 89.81808
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -406.3874269719741
$wsX.Cells.Item(3, 3).Value = 0.09505104843497333
$wsX.Cells.Item(3, 4).Value = 0.6766587302193603
$wsX.Cells.Item(3, 5).Value = 1.9821
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -404.8711513239143
$wsX.Cells.Item(4, 3).Value = 0.09681569305614947
$wsX.Cells.Item(4, 4).Value = 0.9284179101071778
$wsX.Cells.Item(4, 5).Value = 1.56973
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -403.9902935908759
$wsX.Cells.Item(5, 3).Value = 0.06015974628817084
$wsX.Cells.Item(5, 4).Value = 1.5309666000687256
$wsX.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "7" detail ----
$wsX = $wb.Worksheets.Item("7")
$wsX.Cells.Item(2, 4).Value = 0.03952600652770996
$wsX.Cells.Item(2, 5).Value = 83.52166
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -400.1919964156508
$wsX.Cells.Item(3, 3).Value = 0.01569326655826834
$wsX.Cells.Item(3, 4).Value = 3.593181616417114
$wsX.Cells.Item(3, 5).Value = 0.0

# ---- Sheet "8" detail ----
$wsX = $wb.Worksheets.Item("8")
$wsX.Cells.Item(2, 4).Value = 0.04525379928613281
$wsX.Cells.Item(2, 5).Value = 85.27605
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -412.2204295049945
$wsX.Cells.Item(3, 3).Value = 0.05036223552849459
$wsX.Cells.Item(3, 4).Value = 1.5196886937770997
$wsX.Cells.Item(3, 5).Value = 0.59535
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -412.146407489357
$wsX.Cells.Item(4, 3).Value = 0.017960126375595013
$wsX.Cells.Item(4, 4).Value = 3.4196701221275636
$wsX.Cells.Item(4, 5).Value = 0.27881
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -412.09265747622965
$wsX.Cells.Item(5, 3).Value = 0.07472355144900465
$wsX.Cells.Item(5, 4).Value = 4.365206761945313
$wsX.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "9" detail ----
$wsX = $wb.Worksheets.Item("9")
$wsX.Cells.Item(2, 4).Value = 0.04495266033630371
$wsX.Cells.Item(2, 5).Value = 82.46345
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -437.01268675914815
$wsX.Cells.Item(3, 3).Value = 0.00988463101518041
$wsX.Cells.Item(3, 4).Value = 0.1946730905822754
$wsX.Cells.Item(3, 5).Value = 6.51722
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -425.08312443504525
$wsX.Cells.Item(4, 3).Value = 0.0073529907014926515
$wsX.Cells.Item(4, 4).Value = 0.25574909354455566
$wsX.Cells.Item(4, 5).Value = 2.94855
$wsX.Cells.Item(5, 1).Value = 4
$wsX.Cells.Item(5, 2).Value = -414.31946295540973
$wsX.Cells.Item(5, 3).Value = 0.0984312860554829
$wsX.Cells.Item(5, 4).Value = 0.5628914516223145
$wsX.Cells.Item(5, 5).Value = 3.33218
$wsX.Cells.Item(6, 1).Value = 5
$wsX.Cells.Item(6, 2).Value = -408.6154590496442
$wsX.Cells.Item(6, 3).Value = 0.08961188463878143
$wsX.Cells.Item(6, 4).Value = 2.319740486929199
$wsX.Cells.Item(6, 5).Value = 0.55677
$wsX.Cells.Item(7, 1).Value = 6
$wsX.Cells.Item(7, 2).Value = -408.5212229265284
$wsX.Cells.Item(7, 3).Value = 0.0986925316228599
$wsX.Cells.Item(7, 4).Value = 1.6727400641950683
$wsX.Cells.Item(7, 5).Value = 0.74471
$wsX.Cells.Item(8, 1).Value = 7
$wsX.Cells.Item(8, 2).Value = -408.4459826348807
$wsX.Cells.Item(8, 3).Value = 0.04149298119551114
$wsX.Cells.Item(8, 4).Value = 1.7726875184910889
$wsX.Cells.Item(8, 5).Value = 0.0

# ---- Sheet "10" detail ----
$wsX = $wb.Worksheets.Item("10")
$wsX.Cells.Item(2, 4).Value = 0.01759576790344238
$wsX.Cells.Item(2, 5).Value = 84.38256
$wsX.Cells.Item(3, 1).Value = 2
$wsX.Cells.Item(3, 2).Value = -403.581414670978
$wsX.Cells.Item(3, 3).Value = 0.0
$wsX.Cells.Item(3, 4).Value = 5.144143918237793
$wsX.Cells.Item(3, 5).Value = 1.57971
$wsX.Cells.Item(4, 1).Value = 3
$wsX.Cells.Item(4, 2).Value = -403.1039288610207
$wsX.Cells.Item(4, 3).Value = 0.04482826100579968
$wsX.Cells.Item(4, 4).Value = 10.998795292940308
$wsX.Cells.Item(4, 5).Value = 0.0
